$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd drug name "Hydrocone" -> "Hydrocodone" in column A
$ws.Range("A11").Value = "Hydrocodone"

# Update the view state: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("A11").Select()
